$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("2025-08-08")
$src.Copy([System.Reflection.Missing]::Value, $src)
$new = $wb.Worksheets.Item($src.Index + 1)
$new.Name = "2025-08-09"

$data = New-Object 'object[,]' 50,4
$data[0,0] = 1
$data[0,1] = '悪人面したＢ級冒険者 主人公とその幼馴染たちのパパになる'
$data[0,2] = 'こげめ(著者) えんじ(原作) ハラカズヒロ(キャラクター原案)'
$data[0,3] = '第17話-1：「違法奴隷商討伐」'
$data[1,0] = 2
$data[1,1] = '王子様の友達'
$data[1,2] = 'すけろく(著者)'
$data[1,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[2,0] = 3
$data[2,1] = 'ダークサモナーとデキている'
$data[2,2] = '車王(著者)'
$data[2,3] = '第73話'
$data[3,0] = 4
$data[3,1] = '勇者に全部奪われた俺は勇者の母親とパーティを組みました！'
$data[3,2] = '久遠まこと(著者) 石のやっさん(原作)'
$data[3,3] = '第29話'
$data[4,0] = 5
$data[4,1] = 'まんきつしたい常連さん'
$data[4,2] = 'しんみりん(著者)'
$data[4,3] = '第46話後編'
$data[5,0] = 6
$data[5,1] = '淫獄団地'
$data[5,2] = '搾精研究所(原作) 丈山雄為(漫画)'
$data[5,3] = '第49話（後編）'
$data[6,0] = 7
$data[6,1] = 'よくわからないけれど異世界に転生していたようです'
$data[6,2] = '内々けやき あし カオミン'
$data[6,3] = '第137話 よくわからないけれど脱出するみたいです（１）'
$data[7,0] = 8
$data[7,1] = '美人女上司滝沢さん'
$data[7,2] = 'やんBARU(著者)'
$data[7,3] = '第202.5話'
$data[8,0] = 9
$data[8,1] = '治癒魔法の間違った使い方 ~戦場を駆ける回復要員~'
$data[8,2] = '九我山レキ(漫画) くろかた(原作) ＫｅＧ(キャラクター原案)'
$data[8,3] = '第81話その2'
$data[9,0] = 10
$data[9,1] = '解雇された暗黒兵士(30代)のスローなセカンドライフ'
$data[9,2] = '岡沢六十四 るれくちぇ sage・ジョー'
$data[9,3] = '第71話(後編) ダリエルVS.滾り'
$data[10,0] = 11
$data[10,1] = '塔の管理をしてみよう'
$data[10,2] = '盧恩＆雪笠(Friendly Land)(著者) 早秋(原作) 雨神(キャラクター原案)'
$data[10,3] = '第91話後編'
$data[11,0] = 12
$data[11,1] = 'ワンパンマン'
$data[11,2] = '原作/ＯＮＥ 作画/村田雄介'
$data[11,3] = '208撃目'
$data[12,0] = 13
$data[12,1] = '陰キャの俺が席替えでS級美少女に囲まれたら秘密の関係が始まった。'
$data[12,2] = '星野 星野(原作) バラマツヒトミ(漫画) 黒兎 ゆう(キャラクター原案)'
$data[12,3] = '第4話'
$data[13,0] = 14
$data[13,1] = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$data[13,2] = '神原絵理華(漫画) 一森一輝(原作)'
$data[13,3] = '第18話④'
$data[14,0] = 15
$data[14,1] = '異世界ではじめる二拠点生活 ～空間魔法で王都と田舎をいったりきたり～'
$data[14,2] = '丸山りん(漫画) 錬金王(原作) あんべよしろう(キャラクター原案)'
$data[14,3] = 'コミックス第1巻発売告知'
$data[15,0] = 16
$data[15,1] = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$data[15,2] = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$data[15,3] = '第77話 答え'
$data[16,0] = 17
$data[16,1] = 'ベヒモスの花婿'
$data[16,2] = '鈴音ことら(原作) 月並甲介(漫画)'
$data[16,3] = '第8話①'
$data[17,0] = 18
$data[17,1] = 'ヤンデレかと思ったらもっとヤベー女だった'
$data[17,2] = '八木戸マト'
$data[17,3] = '第70話　彼氏を奪われたヤンデレ彼女'
$data[18,0] = 19
$data[18,1] = '彼女にしたい女子一位、の隣で見つけたあまりちゃん'
$data[18,2] = '寝巻ネルゾ(漫画) 裕時悠示(原作) たん旦(キャラクター原案)'
$data[18,3] = 'コミックス1巻のお知らせ'
$data[19,0] = 20
$data[19,1] = '修羅幼女の英雄譚～半端者と言われた傭兵、幼女に転生して成り上がる～'
$data[19,2] = '作画：むらたん 原作：沙城流'
$data[19,3] = '第8話(2)'
$data[20,0] = 21
$data[20,1] = '僕のいけずな婚約者'
$data[20,2] = '冬谷リク(漫画)'
$data[20,3] = '第8話'
$data[21,0] = 22
$data[21,1] = 'ガヴリールドロップアウト'
$data[21,2] = 'うかみ(著者)'
$data[21,3] = '第126話'
$data[22,0] = 23
$data[22,1] = 'センパイ、自宅警備員の雇用はいかがですか？'
$data[22,2] = '漫画：コブラサナギ 原作：二上圭 キャラ原案：日向あずり'
$data[22,3] = '第6話前半'
$data[23,0] = 24
$data[23,1] = '断れない会長は友江くんにだけしてあげたい'
$data[23,2] = '沼地どろまる(著者)'
$data[23,3] = '第29話'
$data[24,0] = 25
$data[24,1] = '序盤で死ぬ最強のサブキャラに転生したので、ゲーム知識で無双する'
$data[24,2] = '作画：マエD 原作：新人'
$data[24,3] = '第5話(3)'
$data[25,0] = 26
$data[25,1] = '最強の少年聖騎士、転生者を狩る'
$data[25,2] = '作画：御塩 原作：宇奈木ユラ'
$data[25,3] = '第6話(3)'
$data[26,0] = 27
$data[26,1] = 'みつばものがたり 呪いの少女と死の輪舞《ロンド》'
$data[26,2] = '堤りん(漫画) 七沢またり(原作) EURA(キャラクター原案)'
$data[26,3] = '第11話：勝利の美酒'
$data[27,0] = 28
$data[27,1] = '佐々木とピーちゃん 異世界でスローライフを楽しもうとしたら、現代で異能バトルに巻き込まれた件 ～魔法少女がアップを始めたようです～'
$data[27,2] = 'ぶんころり(原作) プレジ和尚(漫画) カントク(キャラクター原案)'
$data[27,3] = '第26話①'
$data[28,0] = 29
$data[28,1] = '最強勇者パーティーは愛が知りたい'
$data[28,2] = '山田肌襦袢'
$data[28,3] = '第28話「最後はこぶしがあればいい」'
$data[29,0] = 30
$data[29,1] = '時森さんが無防備です!!'
$data[29,2] = 'たざわ'
$data[29,3] = '第64話'
$data[30,0] = 31
$data[30,1] = 'スキルがなければレベルを上げる～９９がカンストの世界でレベル800万からスタート～'
$data[30,2] = '倉橋ユウス(漫画) 岡沢六十四(原作)'
$data[30,3] = '第52話②'
$data[31,0] = 32
$data[31,1] = 'ダウナーお姉さんは遊びたい'
$data[31,2] = '山鷹景'
$data[31,3] = '第17話'
$data[32,0] = 33
$data[32,1] = 'パワハラギルマスをぶん殴ってブラック聖剣ギルドをクビになったので、辺境で聖剣工房を開くことにした'
$data[32,2] = 'だいたいねむい(原作) まお(漫画)'
$data[32,3] = '第10話①'
$data[33,0] = 34
$data[33,1] = '不老不死少女の苗床旅行記'
$data[33,2] = 'ふじはん(漫画) ルナ・ウサギ(原作)'
$data[33,3] = '第17話前編'
$data[34,0] = 35
$data[34,1] = 'くじ引き特賞：無双ハーレム権'
$data[34,2] = '原作／三木なずな（GA文庫／SBクリエイティブ刊） 漫画／長谷見亮 キャラクター原案／瑠奈璃亜'
$data[34,3] = '第59話-01　777倍の男、魔を統べる者と対峙す！'
$data[35,0] = 36
$data[35,1] = '脱稿するまでオチません'
$data[35,2] = 'ヨシラギ(著者)'
$data[35,3] = '第33話'
$data[36,0] = 37
$data[36,1] = '勇者パーティから追い出された不遇職【罠士】、ユニークスキル【矢印】で最強になる'
$data[36,2] = '作画：たつひこ 原作：白石 有希'
$data[36,3] = '第8話(2)'
$data[37,0] = 38
$data[37,1] = '義妹生活'
$data[37,2] = '三河ごーすと(原作) 奏ユミカ(漫画) Hiten(キャラクター原案)'
$data[37,3] = '第31話-1'
$data[38,0] = 39
$data[38,1] = 'アイツノカノジョ'
$data[38,2] = '肉丸'
$data[38,3] = '7巻発売PR'
$data[39,0] = 40
$data[39,1] = '道にスライムが捨てられていたから連れて帰りました ～おじさんとスライムのほのぼの冒険ライフ～'
$data[39,2] = 'めぐお(漫画) イコ(原作) いもいち(キャラクター原案)'
$data[39,3] = '第2話-1'
$data[40,0] = 41
$data[40,1] = '霜月さんはモブが好き～人見知りな彼女は俺にだけデレ甘い～'
$data[40,2] = '漫画：きぐるみ 原作：八神鏡 キャラクター原案：Roha'
$data[40,3] = '第15話前半'
$data[41,0] = 42
$data[41,1] = '異世界おじさん'
$data[41,2] = '殆ど死んでいる(著者)'
$data[41,3] = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'
$data[42,0] = 43
$data[42,1] = '怠惰な悪辱貴族に転生した俺、シナリオをぶっ壊したら規格外の魔力で最凶になった'
$data[42,2] = '菊池快晴(原作) 小田童馬(作画) 桑島黎音(キャラクター原案)'
$data[42,3] = '第11話'
$data[43,0] = 44
$data[43,1] = 'クロの戦記Ⅱ 異世界転移した僕が最強なのはベッドの上だけのようです'
$data[43,2] = 'サイトウアユム(原作) ユリシロ(漫画) むつみまさと(キャラクター原案)'
$data[43,3] = '第22話-3'
$data[44,0] = 45
$data[44,1] = '神々に育てられしもの、最強となる'
$data[44,2] = '九野十弥(漫画) 羽田遼亮(原作) fame(キャラクター原案)'
$data[44,3] = '第57話'
$data[45,0] = 46
$data[45,1] = 'ライドンキング'
$data[45,2] = '馬場康誌'
$data[45,3] = '第81話 大統領と失われた神器（後編）'
$data[46,0] = 47
$data[46,1] = 'じつは義妹でした。～最近できた義理の弟の距離感がやたら近いわけ～'
$data[46,2] = '堺しょうきち(著者) 白井ムク(原作) 千種みのり(キャラクター原案)'
$data[46,3] = '第36話-2'
$data[47,0] = 48
$data[47,1] = 'チンチンデビルを追え！'
$data[47,2] = 'くぼたふみお'
$data[47,3] = '第３３話　試されしジュノー'
$data[48,0] = 49
$data[48,1] = '最強不敗の神剣使い'
$data[48,2] = '不動らん(漫画) 羽田遼亮(原作) えいひ(キャラクター原案)'
$data[48,3] = '第13話'
$data[49,0] = 50
$data[49,1] = 'りゅうとあまがみ'
$data[49,2] = '角丸柴朗(著者)'
$data[49,3] = '第二話・お肉は何処？①'

$new.Range("A2:D51").Value = $data
